$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Helper: escape plain text for embedding inside a <w:t> element.
# ---------------------------------------------------------------------------
function Escape-Xml {
    param([string]$text)
    return ($text -replace '&', '&amp;' -replace '<', '&lt;' -replace '>', '&gt;')
}

# ---------------------------------------------------------------------------
# Helper: replace the contents of a Range (which must be the exact text
# span of one or more whole runs, WITHOUT the trailing paragraph mark) with
# a freshly built sequence of <w:r> runs, preserving the owning <w:p>'s
# attributes (paraId, rsid, pPr, etc. are untouched because we only ever
# replace the run content, never the paragraph element itself).
# ---------------------------------------------------------------------------
function Replace-RunsXML {
    param(
        [__ComObject]$Range,
        [string]$RunsXml
    )
    $xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p>' + $RunsXml + '</w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $Range.InsertXML($xml)
}

# Build a single plain <w:r> run, optionally with xml:space="preserve" and
# optionally bold.
function Make-Run {
    param(
        [string]$Text,
        [bool]$Preserve = $false,
        [bool]$Bold = $false
    )
    $t = Escape-Xml $Text
    $space = ""
    if ($Preserve) { $space = ' xml:space="preserve"' }
    $rpr = ""
    if ($Bold) { $rpr = "<w:rPr><w:b/><w:bCs/></w:rPr>" }
    return "<w:r>$rpr<w:t$space>$t</w:t></w:r>"
}

# Returns a Range spanning paragraph $index's text, WITHOUT the trailing
# paragraph mark (so InsertXML only rewrites run content, keeping the
# <w:p> element and its attributes intact).
function Get-ParaTextRange {
    param([int]$Index)
    $p = $d.Paragraphs($Index).Range
    return $d.Range($p.Start, $p.End - 1)
}

# ===========================================================================
# 1. Merge the runs in the "Please meet with your team..." paragraph (the
#    spell-check proofErr markers around "Gradescope" disappear along with
#    the run split).
# ===========================================================================
$introText = "Please meet with your team to discuss how you will work together as a team, sign this form (just typing your name is fine), and submit together with the project to Gradescope."
$d.Content.Find.Execute($introText, $true, $false, $false, $false, $false, $true, 1, $false, $introText, 2) | Out-Null

# ===========================================================================
# 2. Add a bold "X" in the middle of the "____" blank for five specific
#    options (turning "____ Label" into "__" + bold "X" + "__ Label").
# ===========================================================================
function Insert-BoldX {
    param([string]$FullText)
    $rest = $FullText.Substring(4)   # text after the four leading underscores
    $rng = $d.Content
    $found = $rng.Find.Execute($FullText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        Write-Output "NOT FOUND: $FullText"
        return
    }
    $whole = $d.Range($rng.Start, $rng.End)
    $runsXml = (Make-Run "__") + (Make-Run "X" $false $true) + (Make-Run ("__" + $rest))
    Replace-RunsXML $whole $runsXml
}

Insert-BoldX "____ Once per week"
Insert-BoldX "____ By deciding in each meeting when the next one will be (and then following up promptly by e-mail in case someone was not present)."
Insert-BoldX "____ By e-mail"
Insert-BoldX "____ By zoom"
Insert-BoldX "____ 24 hours"

# ===========================================================================
# 3. Merge the runs in the "Sign (or type in your name) and date..."
#    paragraph (the spell-check proofErr markers around "Gradescope"
#    disappear here too).
# ===========================================================================
$signText = "Sign (or type in your name) and date the contract below, and keep it on record in a shared space for your team, and have your team leader submit to Gradescope with Milestone 1. "
$d.Content.Find.Execute($signText, $true, $false, $false, $false, $false, $true, 1, $false, $signText, 2) | Out-Null

# ===========================================================================
# 4. Fill in team-member details.
# ===========================================================================

# -- Team Member 1 ----------------------------------------------------------
$rng = Get-ParaTextRange 38   # "Team Member 1 Full name:"
$runsXml = (Make-Run "Team Member 1 Full name:") + (Make-Run " Xiaohai Chen" $true)
Replace-RunsXML $rng $runsXml

$rng = Get-ParaTextRange 39   # "Signed:"
$runsXml = (Make-Run "Signed:") + (Make-Run " " $true) + (Make-Run "Xiaohai Chen")
Replace-RunsXML $rng $runsXml

$rng = Get-ParaTextRange 40   # "Date: "
$runsXml = (Make-Run "Date: " $true) + (Make-Run "3/25/2025")
Replace-RunsXML $rng $runsXml

# -- Team Member 2 ----------------------------------------------------------
$rng = Get-ParaTextRange 42   # "Team Member 2 Full name:"
$runsXml = (Make-Run "Team Member 2 Full name:") + (Make-Run " David Davis" $true)
Replace-RunsXML $rng $runsXml

$rng = Get-ParaTextRange 43   # "Signed:"
$runsXml = (Make-Run "Signed:") + (Make-Run " David Davis" $true)
Replace-RunsXML $rng $runsXml

$rng = Get-ParaTextRange 44   # "Date: "
$runsXml = (Make-Run "Date: " $true) + (Make-Run "3/28/2025")
Replace-RunsXML $rng $runsXml

# ===========================================================================
# 5. Team Member 3's "Date: " + "           " runs collapse into a single
#    run "Date:            ".
# ===========================================================================
$rng = Get-ParaTextRange 48
$runsXml = Make-Run "Date:            " $true
Replace-RunsXML $rng $runsXml

Write-Output "All edits applied"
